# "run prepare & render with final data"
# Updates two of the question-label strings (re-wrapped / re-worded) and
# refreshes the simulated survey-share numbers in the data table with the
# final values produced by the prepare/render pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-wrap the "Governments should actively cooperate..." label (A5) ---
$ws.Range("A5").Value = "`"Governments should actively cooperate to have all countries`nconverge in terms of GDP per capita by the end of the century`""

# --- Replace "Could sign a petition and spread ideas" with the final wording (A6) ---
$ws.Range("A6").Value = "Would support a global movement to tackle CC, tax millionaires,`n and fund LICs (either petition, demonstrate, strike, or donate)"

# --- Refresh numeric results with the final computed data ---
$ws.Range("B2").Value = 0.697581258153402
$ws.Range("K2").Value = 0.686659119892822
$ws.Range("L2").Value = 0.745847844531965
$ws.Range("N2").Value = 0.617013940284116

$ws.Range("B3").Value = 0.641096347070675
$ws.Range("K3").Value = 0.552386382607627
$ws.Range("L3").Value = 0.752151504159109
$ws.Range("N3").Value = 0.575861430623479

$ws.Range("B4").Value = 0.680881448179833
$ws.Range("K4").Value = 0.758076861129753
$ws.Range("L4").Value = 0.688615273248795
$ws.Range("N4").Value = 0.616918649447641

$ws.Range("B5").Value = 0.717914385961719
$ws.Range("K5").Value = 0.702372413171302
$ws.Range("L5").Value = 0.770988593693527
$ws.Range("N5").Value = 0.561039368985046

$ws.Range("B6").Value = 0.675595447215337
$ws.Range("C6").Value = 0.719216740354837
$ws.Range("D6").Value = 0.699222514786681
$ws.Range("E6").Value = 0.688082663981164
$ws.Range("F6").Value = 0.819338712934373
$ws.Range("G6").Value = 0.708420268414952
$ws.Range("H6").Value = 0.74352822863702
$ws.Range("I6").Value = 0.68138828161491
$ws.Range("J6").Value = 0.639531813440066
$ws.Range("K6").Value = 0.557841849059486
$ws.Range("M6").Value = 0.727098526374066
$ws.Range("N6").Value = 0.666651932459956

$ws.Range("B7").Value = 0.681565505605888
$ws.Range("K7").Value = 0.55659974095395
$ws.Range("N7").Value = 0.669595480182533

$ws.Range("B8").Value = 0.456318549602673
$ws.Range("N8").Value = 0.402452424604714

$ws.Range("B9").Value = 0.604334051757766
$ws.Range("K9").Value = 0.584541124690159
$ws.Range("L9").Value = 0.568089649263453
$ws.Range("N9").Value = 0.55410364808233

$ws.Range("B10").Value = 0.55719177445442
